# Counter deal, accept deal and date filter Code updated
# Updates the "Login" sheet's test-data table: the CounterDeal_TC001 row
# (previously row 50) moves down under the Deals_Chat_* rows, and nine new
# rows are appended for CounterDeal_TC002/003, AcceptDeal_TC001/001(2)/002/002(2)
# and DateFilter_TC001/002/003.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- 1. Make room: insert 9 blank rows after the current last row (56) -----
# This mirrors the original edit (new rows appended after row 56), and the
# freshly-inserted rows inherit the row-56 formatting (thin border via style 8)
# the same way Excel's native "Insert" does.
$ws.Rows("57:65").Insert()

# --- 2. Re-write the data rows 50-65 with their final contents -------------
# Column order: A = Automation Test ID, B = UserName/Email, C = Password, D = Expected Result

$rows = @(
  @{ Row=50; A="Deals_Chat_ShipperUser_TC001";      B="rogerdeals21+stan@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=51; A="Deals_Chat_ShipperUser_TC001(2)";    B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=52; A="Deals_Chat_ShipperAdmin_TC002";      B="rogerdeals21+nick@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=53; A="Deals_Chat_ShipperAdmin_TC002(2)";   B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=54; A="Deals_Chat_CarrierUser_TC003";       B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=55; A="Deals_Chat_CarrierUser_TC003(2)";    B="rogerdeals21+stan@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=56; A="CounterDeal_TC001";                  B="rogerdeals21+stan@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=57; A="CounterDeal_TC002";                  B="rogerdeals21+rick@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=58; A="CounterDeal_TC003";                  B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=59; A="AcceptDeal_TC001";                   B="rogerdeals21+stan@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=60; A="AcceptDeal_TC001(2)";                B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=61; A="AcceptDeal_TC002(2)";                B="rogerdeals21+stan@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=62; A="AcceptDeal_TC002";                   B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=63; A="DateFilter_TC001";                   B="rogerdeals21+stan@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=64; A="DateFilter_TC002";                   B="rogerdeals21+rick@gmail.com"; C="arewethere?"; D="Login successful" },
  @{ Row=65; A="DateFilter_TC003";                   B="rogerdeals21+john@gmail.com"; C="arewethere?"; D="Login successful" }
)

foreach ($r in $rows) {
  $ws.Range("A$($r.Row)").Value = $r.A
  $ws.Range("B$($r.Row)").Value = $r.B
  $ws.Range("C$($r.Row)").Value = $r.C
  $ws.Range("D$($r.Row)").Value = $r.D
}

# --- 3. Fix up the vertical-centred style (style index 9) on column A ------
# Rows 50-51 keep it, 52-55 lose it, and 56-65 (the CounterDeal/AcceptDeal/
# DateFilter block) gain it, matching the source formatting used elsewhere
# in this test-data column.
$ws.Range("A53").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("A55").PasteSpecial(-4122)

$ws.Range("A50").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("A58").PasteSpecial(-4122)
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A65").PasteSpecial(-4122)

# --- 4. View state: scroll position & selection match the saved workbook ---
$ws.Activate()
$ws.Range("A50:A55").Select()
$excel.ActiveWindow.ScrollRow = 46
